$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H4: 1014 -> 507
$ws.Range("H4").Value = 507

# E6: "Sweat Beast" -> "Nyakuza Mask`nA hat in time"
$ws.Range("E6").Value = "Nyakuza Mask`nA hat in time"

# H6: discord link -> empty
$ws.Range("H6").Value = ""

# H8: empty -> dunkbin link
$ws.Range("H8").Value = "https://dunkbin.com/img/623.png"

# H9: discord link -> empty
$ws.Range("H9").Value = ""

# D12: clown_noes666 -> snekiecr8
$ws.Range("D12").Value = "snekiecr8"

# E12: 469660616 -> 28059068
$ws.Range("E12").Value = 28059068

# F12: 118 days -> 616 days
$ws.Range("F12").Value = "616 days"

# D14: clown noes666 -> snekie
$ws.Range("D14").Value = "snekie"

# D15: Clown_noes666 -> Snekiecr8
$ws.Range("D15").Value = "Snekiecr8"

# K16: 1014 -> 1027
$ws.Range("K16").Value = 1027
